$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 1069
$ws.Range("I12").Value = 1465
$ws.Range("J12").Value = 475
$ws.Range("K12").Value = 1465
$ws.Range("L12").Value = 475
$ws.Range("M12").Value = -1295
$ws.Range("N12").Value = -815
$ws.Range("H13").Value = 1013.5
$ws.Range("I13").Value = 940
$ws.Range("J13").Value = 1087
$ws.Range("K13").Value = 940
$ws.Range("L13").Value = 1087
$ws.Range("M13").Value = -771
$ws.Range("N13").Value = -1425
$ws.Range("H40").Value = 5760
$ws.Range("J40").Value = 6277.778
$ws.Range("L40").Value = 6277.778
$ws.Range("N40").Value = -6627.778
$ws.Range("H96").Value = 1474
$ws.Range("I96").Value = 2335
$ws.Range("J96").Value = 900
$ws.Range("K96").Value = 7005
$ws.Range("L96").Value = 2700
$ws.Range("M96").Value = -5632
$ws.Range("N96").Value = -5446
$ws.Range("H113").Value = 3049.25
$ws.Range("J113").Value = 3300
$ws.Range("L113").Value = 3300
$ws.Range("N113").Value = -9808
$ws.Range("H135").Value = 250
$ws.Range("I135").Value = 250
$ws.Range("K135").Value = 2250
$ws.Range("M135").Value = 285

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2213
$ws.Range("I32").Value = 2289.4
$ws.Range("K32").Value = 2289.4
$ws.Range("M32").Value = -2002.4
$ws.Range("H102").Value = 0
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("M102").ClearContents()
$ws.Range("N102").Value = 0
$ws.Range("H122").Value = 2560.3635
$ws.Range("I122").Value = 2541
$ws.Range("J122").Value = 2571.4285
$ws.Range("K122").Value = 7623
$ws.Range("L122").Value = 7714.2855
$ws.Range("M122").Value = -5173
$ws.Range("N122").Value = -12614.2855

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 765.1111
$ws.Range("I86").Value = 721.2
$ws.Range("J86").Value = 820
$ws.Range("K86").Value = 721.2
$ws.Range("L86").Value = 820
$ws.Range("M86").Value = 401.8
$ws.Range("N86").Value = -3066
$ws.Range("H89").Value = 765.1111
$ws.Range("I89").Value = 721.2
$ws.Range("J89").Value = 820
$ws.Range("K89").Value = 3606
$ws.Range("L89").Value = 4100
$ws.Range("M89").Value = 2010
$ws.Range("N89").Value = -15332
$ws.Range("H94").Value = 8684.143
$ws.Range("I94").Value = 7747.5
$ws.Range("K94").Value = 7747.5
$ws.Range("M94").Value = -7296.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8209.056
$ws.Range("I31").Value = 8930.429
$ws.Range("K31").Value = 8930.429
$ws.Range("M31").Value = -8635.429
$ws.Range("H34").Value = 8209.056
$ws.Range("I34").Value = 8930.429
$ws.Range("K34").Value = 8930.429
$ws.Range("M34").Value = -8728.429
$ws.Range("H105").Value = 1224.75
$ws.Range("I105").Value = 1166
$ws.Range("J105").Value = 1260
$ws.Range("K105").Value = 1166
$ws.Range("L105").Value = 1260
$ws.Range("M105").Value = 581
$ws.Range("N105").Value = -4754
$ws.Range("H122").Value = 1515.6666
$ws.Range("I122").Value = 1273.5
$ws.Range("K122").Value = 3820.5
$ws.Range("M122").Value = -1370.5
$ws.Range("H134").Value = 5207.25
$ws.Range("I134").Value = 5207.25
$ws.Range("K134").Value = 15621.75
$ws.Range("M134").Value = -13086.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 195
$ws.Range("J34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("N34").ClearContents()
$ws.Range("H39").Value = 5000
$ws.Range("J39").Value = 5000
$ws.Range("L39").Value = 15000
$ws.Range("N39").Value = -15588
$ws.Range("H45").Value = 150
$ws.Range("I45").Value = 150
$ws.Range("K45").Value = 450
$ws.Range("M45").Value = 82
$ws.Range("H129").Value = 1524.1428
$ws.Range("J129").Value = 1734
$ws.Range("L129").Value = 5202
$ws.Range("N129").Value = -15202

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 5284.857
$ws.Range("I126").Value = 5332.3335
$ws.Range("J126").Value = 5000
$ws.Range("K126").Value = 15997.0005
$ws.Range("L126").Value = 15000
$ws.Range("M126").Value = -13527.0005
$ws.Range("N126").Value = -19940
$ws.Range("H132").Value = 4506.875
$ws.Range("I132").Value = 4833.3335
$ws.Range("K132").Value = 14500.0005
$ws.Range("M132").Value = -11970.0005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1901
$ws.Range("I46").Value = 2000
$ws.Range("J46").Value = 1802
$ws.Range("K46").Value = 2000
$ws.Range("L46").Value = 1802
$ws.Range("M46").Value = -1812
$ws.Range("N46").Value = -2178
$ws.Range("H122").Value = 5051
$ws.Range("I122").Value = 5401.3335
$ws.Range("K122").Value = 16204.0005
$ws.Range("M122").Value = -13754.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 1074009
$ws.Range("I2").Value = 1569479.2
$ws.Range("J2").Value = 490.33334
$ws.Range("K2").Value = 1569479.2
$ws.Range("L2").Value = 490.33334
$ws.Range("M2").Value = -1569367.2
$ws.Range("N2").Value = -714.33334
$ws.Range("H4").Value = 20725.5
$ws.Range("I4").Value = 28985.5
$ws.Range("K4").Value = 28985.5
$ws.Range("M4").Value = -28872.5
$ws.Range("H62").Value = 3498
$ws.Range("I62").Value = 3498
$ws.Range("K62").Value = 3498
$ws.Range("M62").Value = -2874
$ws.Range("H65").Value = 3498
$ws.Range("I65").Value = 3498
$ws.Range("K65").Value = 17490
$ws.Range("M65").Value = -14370
$ws.Range("H100").Value = 1699.2
$ws.Range("I100").Value = 499
$ws.Range("J100").Value = 3499.5
$ws.Range("K100").Value = 998
$ws.Range("L100").Value = 6999
$ws.Range("M100").Value = -457
$ws.Range("N100").Value = -8081
$ws.Range("H132").Value = 2350.818
$ws.Range("I132").Value = 2801.75
$ws.Range("J132").Value = 1148.3334
$ws.Range("K132").Value = 8405.25
$ws.Range("L132").Value = 3445.0002
$ws.Range("M132").Value = -5875.25
$ws.Range("N132").Value = -8505.0002
$ws.Range("H136").Value = 1071
$ws.Range("I136").Value = 1071
$ws.Range("K136").Value = 3213
$ws.Range("M136").Value = -663
